$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 value from 4.2 to 6.2
$ws.Range("C2").Value = 6.2

# Update selection to C5
$ws.Range("C5").Select()

# New rows data: column A is a name (shared string), column B is a numeric value
$names = @(
    "אן מרש",
    "תאיו ורד",
    "יהלי דוייב",
    "אורי שטרנברג",
    "אורי שטרנברג",
    "אורי שטרנברג",
    "רומי הרשקוביץ",
    "אביב ואסקז",
    "הילה שולויס",
    "קרן רינת פביאן",
    "ליהי בראל",
    "ליאם דיין",
    "לינוי קוסטיקה",
    "מעיין סטרוזר",
    "יולי קזמה",
    "שלו דיין",
    "איתי הראל",
    "מעיין סטרוזר",
    "מעיין סטרוזר"
)

$values = @(1,1,1,1,6,6,1,1,1,1,1,1,1,1,1,1,1,6,6)

$startRow = 496
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
